$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the D column tail (D32:D34) - these values are no longer needed
$ws.Range("D32:D34").ClearContents()

# Clear out rows 52-58 entirely (B and C columns) - fully empty rows drop from the sheet
$ws.Range("B52:C58").ClearContents()

# Clear the values in the last three rows (B59:C61) but keep their existing style
$ws.Range("B59:C61").ClearContents()

# Update the selection to match the new active cell / scroll position
$ws.Range("D32").Select()
